# Append job-listing rows scraped at 2025-11-08 01:15:48 JST.
# Existing rows shift/refresh and 2 brand-new postings are added,
# while one stale test row ("初回 ssss") from the previous run is gone.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out any existing hyperlink objects first so we can rebuild them
# cleanly against the refreshed row layout (avoids stale row associations).
$ws.Hyperlinks.Delete()

# Row 2: 5429252
$ws.Range("A2").Value = '2025-11-08 01:15:48'
$ws.Range("B2").Value = '中古ブランド品リサーチとEC出品作業をAIで自動化するツール開発'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5429252'
$ws.Range("G2").Value = 480
$ws.Range("H2").Value = '🔥AI,Ai ◆ツール,開発'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5429252') | Out-Null
$ws.Range("F2").Style = "Hyperlink"

# Row 3: 5428695
$ws.Range("A3").Value = '2025-11-08 01:15:48'
$ws.Range("B3").Value = '専門データ分析:AIコスト最適化設計と厳格な機密保持を必須とするWebシステム開発(段階的継続発注)'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5428695'
$ws.Range("G3").Value = 403
$ws.Range("H3").Value = '🔥AI,Ai ◆開発,システム開発'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5428695') | Out-Null
$ws.Range("F3").Style = "Hyperlink"

# Row 4: 5429181
$ws.Range("A4").Value = '2025-11-08 01:15:48'
$ws.Range("B4").Value = '初回 APIを利用したPowerShellもしくはPythonのスクリプト作成'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5429181'
$ws.Range("G4").Value = 380
$ws.Range("H4").Value = '🔥Python,API'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5429181') | Out-Null
$ws.Range("F4").Style = "Hyperlink"

# Row 5: 5429304
$ws.Range("A5").Value = '2025-11-08 01:15:48'
$ws.Range("B5").Value = '複数の見積書から情報抜き出して別表に書き出す作業の自動化 Excel VBAツール開発依頼'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5429304'
$ws.Range("G5").Value = 208
$ws.Range("H5").Value = '◆ツール,開発'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5429304') | Out-Null
$ws.Range("F5").Style = "Hyperlink"

# Row 6: 5428871
$ws.Range("A6").Value = '2025-11-08 01:15:48'
$ws.Range("B6").Value = '【急募】ECサイト管理システムの開発・保守業務依頼'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5428871'
$ws.Range("G6").Value = 143
$ws.Range("H6").Value = '◆開発 ◇サイト'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5428871') | Out-Null
$ws.Range("F6").Style = "Hyperlink"

# Row 7: 5429220
$ws.Range("A7").Value = '2025-11-08 01:15:48'
$ws.Range("B7").Value = '【急募】モバイルアプリ テスト業務 委託募集(3 - 4週間)'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5429220'
$ws.Range("G7").Value = 38
$ws.Range("H7").Value = '◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5429220') | Out-Null
$ws.Range("F7").Style = "Hyperlink"

# Row 8: 5429495
$ws.Range("A8").Value = '2025-11-08 01:15:48'
$ws.Range("B8").Value = '【急募】既存で作成済みのAccessデータベースをPower Apps等のアプリに移行したい'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5429495'
$ws.Range("G8").Value = 33
$ws.Range("H8").Value = '◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5429495') | Out-Null
$ws.Range("F8").Style = "Hyperlink"

# Row 9: 5429159
$ws.Range("A9").Value = '2025-11-08 01:15:48'
$ws.Range("B9").Value = '【急募】WordPressにe-SCOTT決済機能を導入'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5429159'
$ws.Range("G9").Value = 33
$ws.Range("H9").Value = '○WordPress'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5429159') | Out-Null
$ws.Range("F9").Style = "Hyperlink"

# Row 10: 5429157
$ws.Range("A10").Value = '2025-11-08 01:15:48'
$ws.Range("B10").Value = '【急募】WordPressにe-SCOTT決済機能を導入'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5429157'
$ws.Range("G10").Value = 33
$ws.Range("H10").Value = '○WordPress'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5429157') | Out-Null
$ws.Range("F10").Style = "Hyperlink"

# Row 11: 5429335
$ws.Range("A11").Value = '2025-11-08 01:15:48'
$ws.Range("B11").Value = '【フルスタックエンジニア】 働きながらスキルアップもできるEC業界で活躍してみませんか?'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5429335'
$ws.Range("G11").Value = 25
$ws.Range("H11").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5429335') | Out-Null
$ws.Range("F11").Style = "Hyperlink"

# Row 12: 5428756
$ws.Range("A12").Value = '2025-11-08 01:15:48'
$ws.Range("B12").Value = '【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5428756'
$ws.Range("G12").Value = 25
$ws.Range("H12").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5428756') | Out-Null
$ws.Range("F12").Style = "Hyperlink"

# Row 13: 5428755
$ws.Range("A13").Value = '2025-11-08 01:15:48'
$ws.Range("B13").Value = '【リーダー募集×リモートOK】SRE/インフラエンジニア(Google Cloud/長期金融系案件)'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5428755'
$ws.Range("G13").Value = 25
$ws.Range("H13").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5428755') | Out-Null
$ws.Range("F13").Style = "Hyperlink"

# Row 14: 5428970
$ws.Range("A14").Value = '2025-11-08 01:15:48'
$ws.Range("B14").Value = '【急募】Wordプレスエラー修正のプロを探しています!'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '~ 5,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5428970'
$ws.Range("G14").Value = 10
$ws.Range("H14").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5428970') | Out-Null
$ws.Range("F14").Style = "Hyperlink"

# Row 15: 5428509
$ws.Range("A15").Value = '2025-11-08 01:15:48'
$ws.Range("B15").Value = '【急募】Googleworkスペース・ハブスポットのサーバー設定依頼'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5428509'
$ws.Range("G15").Value = 10
$ws.Range("H15").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5428509') | Out-Null
$ws.Range("F15").Style = "Hyperlink"
